$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44446
$ws.Cells.Item(2, 10).Value = 160
$ws.Cells.Item(2, 11).Value = 12500
$ws.Cells.Item(2, 12).Value = 13000
$ws.Cells.Item(2, 13).Value = 12750
$ws.Cells.Item(2, 16).Value = 319
$ws.Cells.Item(3, 4).Value = 44475
$ws.Cells.Item(3, 11).Value = 11000
$ws.Cells.Item(3, 12).Value = 12000
$ws.Cells.Item(3, 13).Value = 11500
$ws.Cells.Item(3, 16).Value = 288
$ws.Cells.Item(4, 4).Value = 44488
$ws.Cells.Item(4, 8).Value = 'Madrigal'
$ws.Cells.Item(4, 9).Value = 'Primera'
$ws.Cells.Item(4, 10).Value = 100
$ws.Cells.Item(4, 11).Value = 11000
$ws.Cells.Item(4, 12).Value = 12000
$ws.Cells.Item(4, 13).Value = 11500
$ws.Cells.Item(4, 16).Value = 288
$ws.Cells.Item(5, 4).Value = 44453
$ws.Cells.Item(5, 10).Value = 160
$ws.Cells.Item(5, 13).Value = 12750
$ws.Cells.Item(5, 16).Value = 319
$ws.Cells.Item(6, 4).Value = 44516
$ws.Cells.Item(7, 4).Value = 44515
$ws.Cells.Item(7, 11).Value = 11000
$ws.Cells.Item(7, 12).Value = 12000
$ws.Cells.Item(7, 13).Value = 11500
$ws.Cells.Item(7, 16).Value = 288
$ws.Cells.Item(8, 4).Value = 44435
$ws.Cells.Item(8, 10).Value = 120
$ws.Cells.Item(8, 11).Value = 14000
$ws.Cells.Item(8, 12).Value = 15000
$ws.Cells.Item(8, 13).Value = 14500
$ws.Cells.Item(8, 16).Value = 362
$ws.Cells.Item(9, 4).Value = 44468
$ws.Cells.Item(9, 10).Value = 60
$ws.Cells.Item(9, 11).Value = 12000
$ws.Cells.Item(9, 12).Value = 13000
$ws.Cells.Item(9, 13).Value = 12500
$ws.Cells.Item(9, 16).Value = 312
$ws.Cells.Item(10, 4).Value = 44484
$ws.Cells.Item(10, 10).Value = 120
$ws.Cells.Item(11, 4).Value = 44425
$ws.Cells.Item(11, 10).Value = 120
$ws.Cells.Item(11, 11).Value = 14000
$ws.Cells.Item(11, 12).Value = 15000
$ws.Cells.Item(11, 13).Value = 14500
$ws.Cells.Item(11, 15).Value = 'Región del Maule'
$ws.Cells.Item(11, 16).Value = 362
$ws.Cells.Item(12, 4).Value = 44417
$ws.Cells.Item(12, 11).Value = 15000
$ws.Cells.Item(12, 12).Value = 16000
$ws.Cells.Item(12, 13).Value = 15500
$ws.Cells.Item(12, 16).Value = 388
$ws.Cells.Item(13, 4).Value = 44490
$ws.Cells.Item(13, 10).Value = 100
$ws.Cells.Item(14, 4).Value = 44495
$ws.Cells.Item(14, 11).Value = 11000
$ws.Cells.Item(14, 12).Value = 12000
$ws.Cells.Item(14, 13).Value = 11500
$ws.Cells.Item(14, 16).Value = 288
$ws.Cells.Item(15, 4).Value = 44503
$ws.Cells.Item(15, 10).Value = 160
$ws.Cells.Item(16, 4).Value = 44498
$ws.Cells.Item(16, 10).Value = 60
$ws.Cells.Item(16, 11).Value = 10500
$ws.Cells.Item(16, 12).Value = 11000
$ws.Cells.Item(16, 13).Value = 10750
$ws.Cells.Item(16, 16).Value = 269
$ws.Cells.Item(17, 4).Value = 44420
$ws.Cells.Item(17, 11).Value = 13000
$ws.Cells.Item(17, 12).Value = 14000
$ws.Cells.Item(17, 13).Value = 13500
$ws.Cells.Item(17, 16).Value = 338
$ws.Cells.Item(18, 4).Value = 44508
$ws.Cells.Item(18, 10).Value = 160
$ws.Cells.Item(19, 4).Value = 44467
$ws.Cells.Item(19, 10).Value = 160
$ws.Cells.Item(19, 13).Value = 11500
$ws.Cells.Item(19, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(19, 16).Value = 288
$ws.Cells.Item(20, 4).Value = 44510
$ws.Cells.Item(20, 11).Value = 11000
$ws.Cells.Item(20, 12).Value = 12000
$ws.Cells.Item(20, 13).Value = 11500
$ws.Cells.Item(20, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(20, 16).Value = 288
$ws.Cells.Item(21, 4).Value = 44494
$ws.Cells.Item(21, 10).Value = 120
$ws.Cells.Item(22, 4).Value = 44426
$ws.Cells.Item(22, 13).Value = 13500
$ws.Cells.Item(22, 15).Value = 'Región del Maule'
$ws.Cells.Item(22, 16).Value = 338
$ws.Cells.Item(23, 4).Value = 44473
$ws.Cells.Item(23, 10).Value = 160
$ws.Cells.Item(24, 4).Value = 44512
$ws.Cells.Item(24, 11).Value = 11000
$ws.Cells.Item(24, 12).Value = 12000
$ws.Cells.Item(24, 13).Value = 11500
$ws.Cells.Item(24, 16).Value = 288
$ws.Cells.Item(26, 4).Value = 44455
$ws.Cells.Item(26, 10).Value = 100
$ws.Cells.Item(26, 11).Value = 13000
$ws.Cells.Item(26, 12).Value = 14000
$ws.Cells.Item(26, 13).Value = 13500
$ws.Cells.Item(26, 16).Value = 338
$ws.Cells.Item(27, 4).Value = 44399
$ws.Cells.Item(27, 8).Value = 'Española'
$ws.Cells.Item(27, 9).Value = 'Segunda'
$ws.Cells.Item(27, 10).Value = 120
$ws.Cells.Item(27, 11).Value = 15500
$ws.Cells.Item(27, 12).Value = 16000
$ws.Cells.Item(27, 13).Value = 15750
$ws.Cells.Item(27, 16).Value = 394
$ws.Cells.Item(28, 4).Value = 44487
$ws.Cells.Item(28, 10).Value = 100
$ws.Cells.Item(28, 11).Value = 11000
$ws.Cells.Item(28, 12).Value = 12000
$ws.Cells.Item(28, 13).Value = 11500
$ws.Cells.Item(28, 16).Value = 288
$ws.Cells.Item(29, 4).Value = 44505
$ws.Cells.Item(30, 4).Value = 44432
$ws.Cells.Item(30, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(31, 4).Value = 44496
$ws.Cells.Item(31, 10).Value = 120
$ws.Cells.Item(31, 11).Value = 11000
$ws.Cells.Item(31, 12).Value = 12000
$ws.Cells.Item(31, 13).Value = 11500
$ws.Cells.Item(31, 16).Value = 288
$ws.Cells.Item(32, 4).Value = 44427
$ws.Cells.Item(32, 10).Value = 120
$ws.Cells.Item(32, 11).Value = 13000
$ws.Cells.Item(32, 12).Value = 14000
$ws.Cells.Item(32, 13).Value = 13500
$ws.Cells.Item(32, 16).Value = 338
$ws.Cells.Item(33, 4).Value = 44491
$ws.Cells.Item(33, 10).Value = 100
$ws.Cells.Item(33, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(35, 4).Value = 44454
$ws.Cells.Item(35, 11).Value = 13000
$ws.Cells.Item(35, 12).Value = 14000
$ws.Cells.Item(35, 13).Value = 13500
$ws.Cells.Item(35, 16).Value = 338
